$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.2808873333333333
$ws.Range("H2").Value = 0.842662
$ws.Range("I2").Value = 0.5595554696739399
$ws.Range("J2").Value = 0.5595554696739399
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.5203476666666668
$ws.Range("N2").Value = 1.561043
$ws.Range("O2").Value = 0.004105934376266647
$ws.Range("P2").Value = 0.004105934376266647
$ws.Range("Q2").Value = 0.1461590684962222
$ws.Range("R2").Value = 1.315431616466
$ws.Range("S2").Value = 0.002297498038362259
$ws.Range("T2").Value = 0.002297498038362259

$ws.Range("G3").Value = 0.2808873333333333
$ws.Range("H3").Value = 0.842662
$ws.Range("I3").Value = 0.5595554696739399
$ws.Range("J3").Value = 0.5595554696739399
$ws.Range("O3").Value = 0.8361295370252257
$ws.Range("P3").Value = 0.8361295370252259
$ws.Range("Q3").Value = 29.76372807616644
$ws.Range("R3").Value = 267.873552685498
$ws.Range("S3").Value = 0.4678608557984041
$ws.Range("T3").Value = 0.4678608557984041

$ws.Range("G4").Value = 0.2808873333333333
$ws.Range("H4").Value = 0.842662
$ws.Range("I4").Value = 0.5595554696739399
$ws.Range("J4").Value = 0.5595554696739399
$ws.Range("M4").Value = 20.24706
$ws.Range("N4").Value = 60.74118
$ws.Range("O4").Value = 0.1597645285985076
$ws.Range("P4").Value = 0.1597645285985076
$ws.Range("Q4").Value = 5.68714269124
$ws.Range("R4").Value = 51.18428422116
$ws.Range("S4").Value = 0.08939711583717351
$ws.Range("T4").Value = 0.08939711583717354

$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.2210956666666667
$ws.Range("H5").Value = 0.663287
$ws.Range("I5").Value = 0.4404445303260602
$ws.Range("J5").Value = 0.4404445303260602
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5203476666666668
$ws.Range("N5").Value = 1.561043
$ws.Range("O5").Value = 0.004105934376266647
$ws.Range("P5").Value = 0.004105934376266647
$ws.Range("Q5").Value = 0.1150466142601111
$ws.Range("R5").Value = 1.035419528341
$ws.Range("S5").Value = 0.001808436337904388
$ws.Range("T5").Value = 0.001808436337904388

$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.2210956666666667
$ws.Range("H6").Value = 0.663287
$ws.Range("I6").Value = 0.4404445303260602
$ws.Range("J6").Value = 0.4404445303260602
$ws.Range("O6").Value = 0.8361295370252257
$ws.Range("P6").Value = 0.8361295370252259
$ws.Range("Q6").Value = 23.42801016831922
$ws.Range("R6").Value = 210.852091514873
$ws.Range("S6").Value = 0.3682686812268218
$ws.Range("T6").Value = 0.3682686812268218

$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.2210956666666667
$ws.Range("H7").Value = 0.663287
$ws.Range("I7").Value = 0.4404445303260602
$ws.Range("J7").Value = 0.4404445303260602
$ws.Range("M7").Value = 20.24706
$ws.Range("N7").Value = 60.74118
$ws.Range("O7").Value = 0.1597645285985076
$ws.Range("P7").Value = 0.1597645285985076
$ws.Range("Q7").Value = 4.47653722874
$ws.Range("R7").Value = 40.28883505866
$ws.Range("S7").Value = 0.0703674127613341
$ws.Range("T7").Value = 0.0703674127613341
